$d = $word.ActiveDocument

$d.Content.Find.Execute("998÷9=110, 8", $true, $false, $false, $false, $false, $true, 1, $false, "471÷4=117, 3", 2) | Out-Null
$d.Content.Find.Execute("905÷3=301, 2", $true, $false, $false, $false, $false, $true, 1, $false, "681÷8=85, 1", 2) | Out-Null
$d.Content.Find.Execute("404÷9=44, 8", $true, $false, $false, $false, $false, $true, 1, $false, "898÷7=128, 2", 2) | Out-Null
$d.Content.Find.Execute("742÷2=371, 0", $true, $false, $false, $false, $false, $true, 1, $false, "891÷6=148, 3", 2) | Out-Null
$d.Content.Find.Execute("723÷4=180, 3", $true, $false, $false, $false, $false, $true, 1, $false, "677÷5=135, 2", 2) | Out-Null
$d.Content.Find.Execute("253÷4=63, 1", $true, $false, $false, $false, $false, $true, 1, $false, "698÷9=77, 5", 2) | Out-Null
$d.Content.Find.Execute("491÷8=61, 3", $true, $false, $false, $false, $false, $true, 1, $false, "382÷6=63, 4", 2) | Out-Null
$d.Content.Find.Execute("809÷8=101, 1", $true, $false, $false, $false, $false, $true, 1, $false, "938÷2=469, 0", 2) | Out-Null
$d.Content.Find.Execute("158÷8=19, 6", $true, $false, $false, $false, $false, $true, 1, $false, "739÷7=105, 4", 2) | Out-Null
$d.Content.Find.Execute("316÷6=52, 4", $true, $false, $false, $false, $false, $true, 1, $false, "124÷8=15, 4", 2) | Out-Null
$d.Content.Find.Execute("233÷9=25, 8", $true, $false, $false, $false, $false, $true, 1, $false, "354÷8=44, 2", 2) | Out-Null
$d.Content.Find.Execute("120÷7=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "307÷5=61, 2", 2) | Out-Null
$d.Content.Find.Execute("730÷2=365, 0", $true, $false, $false, $false, $false, $true, 1, $false, "351÷5=70, 1", 2) | Out-Null
$d.Content.Find.Execute("815÷4=203, 3", $true, $false, $false, $false, $false, $true, 1, $false, "944÷5=188, 4", 2) | Out-Null
$d.Content.Find.Execute("949÷3=316, 1", $true, $false, $false, $false, $false, $true, 1, $false, "789÷4=197, 1", 2) | Out-Null
$d.Content.Find.Execute("557÷2=278, 1", $true, $false, $false, $false, $false, $true, 1, $false, "395÷4=98, 3", 2) | Out-Null
$d.Content.Find.Execute("197÷8=24, 5", $true, $false, $false, $false, $false, $true, 1, $false, "818÷4=204, 2", 2) | Out-Null
$d.Content.Find.Execute("260÷5=52, 0", $true, $false, $false, $false, $false, $true, 1, $false, "696÷2=348, 0", 2) | Out-Null
$d.Content.Find.Execute("512÷4=128, 0", $true, $false, $false, $false, $false, $true, 1, $false, "327÷5=65, 2", 2) | Out-Null
$d.Content.Find.Execute("523÷7=74, 5", $true, $false, $false, $false, $false, $true, 1, $false, "448÷9=49, 7", 2) | Out-Null
$d.Content.Find.Execute("863÷6=143, 5", $true, $false, $false, $false, $false, $true, 1, $false, "758÷5=151, 3", 2) | Out-Null
$d.Content.Find.Execute("266÷3=88, 2", $true, $false, $false, $false, $false, $true, 1, $false, "800÷3=266, 2", 2) | Out-Null
$d.Content.Find.Execute("977÷4=244, 1", $true, $false, $false, $false, $false, $true, 1, $false, "288÷2=144, 0", 2) | Out-Null
$d.Content.Find.Execute("634÷7=90, 4", $true, $false, $false, $false, $false, $true, 1, $false, "760÷2=380, 0", 2) | Out-Null
$d.Content.Find.Execute("433÷3=144, 1", $true, $false, $false, $false, $false, $true, 1, $false, "387÷8=48, 3", 2) | Out-Null
